$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy style from existing header cell (e.g. E1) so it matches s="1"
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Per-row time_taken values (plain, unstyled cells matching the diff)
$times = @(
    "2021-10-05 10:50:13.356639",
    "2021-10-05 10:50:13.356653",
    "2021-10-05 10:50:13.356657",
    "2021-10-05 10:50:13.356660",
    "2021-10-05 10:50:13.356663",
    "2021-10-05 10:50:13.356667",
    "2021-10-05 10:50:13.356670",
    "2021-10-05 10:50:13.356673",
    "2021-10-05 10:50:13.356676",
    "2021-10-05 10:50:13.356679",
    "2021-10-05 10:50:13.356682",
    "2021-10-05 10:50:13.356685",
    "2021-10-05 10:50:13.356688",
    "2021-10-05 10:50:13.356691"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
